$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the completed "Boost Bootstrap script should copy BJam to
#    HadesMem root." list item, and move the Word-managed "_GoBack" bookmark
#    (which tracks the last edit point) from its old location to the end of
#    the paragraph that now becomes the last item before "New Modules".
# ---------------------------------------------------------------------------

# The _GoBack bookmark currently sits after "E.g. " earlier in the document;
# remove it from there first (Word will recreate it at the new edit point).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate and delete the whole completed-item paragraph (including its own
# paragraph mark), which merges it away and leaves the previous bullet
# ("...Iterators, Scanner, PeLib, FindPattern, etc).") as the last one.
$targetRange = $d.Content.Duplicate
$targetRange.Find.Execute(
    "Boost Bootstrap script should copy BJam to HadesMem root.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$targetPara = $targetRange.Paragraphs(1).Range
$targetPara.Delete()

# Re-anchor the "_GoBack" bookmark at the end of the preceding paragraph's
# text (i.e. right after the final "." and before the paragraph mark).
$prevRange = $d.Content.Duplicate
$prevRange.Find.Execute(
    "Performance improvements in potential bottlenecks",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$prevPara = $prevRange.Paragraphs(1).Range
$insertPos = $prevPara.End - 1

# A zero-length range placed exactly at that boundary can't be bookmarked
# directly, so briefly insert a marker character after the boundary, anchor
# the (now non-boundary) zero-length bookmark before it, then remove the
# marker again - leaving a clean, empty "_GoBack" bookmark in place.
$markerRange = $d.Range($insertPos, $insertPos)
$markerRange.InsertAfter("X")
$bookmarkRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
$d.Range($insertPos, $insertPos + 1).Delete()

# ---------------------------------------------------------------------------
# 2. Touch the footnote/endnote machinery so the (until now absent)
#    word/footnotes.xml and word/endnotes.xml parts get created with their
#    standard separator / continuation-separator boilerplate, matching a
#    full-fidelity Word save.
# ---------------------------------------------------------------------------
$fnAnchor = $d.Range(0, 0)
$tempFootnote = $d.Footnotes.Add($fnAnchor, "", "x")
$tempFootnote.Delete()

# ---------------------------------------------------------------------------
# 3. Touch the built-in Header/Footer paragraph styles (and their linked
#    "... Char" character styles) so they get written into styles.xml, again
#    matching what a full-fidelity Word save produces.
# ---------------------------------------------------------------------------
$styleAnchorPara = $d.Paragraphs(1)
$styleAnchorOriginal = $styleAnchorPara.Range.Style

$styleAnchorPara.Range.Style = "Header"
$styleAnchorPara.Range.Style = $styleAnchorOriginal

$headerChar = $d.Styles.Add("HeaderChar", 2)
$headerChar.NameLocal = "Header Char"
$headerChar.BaseStyle = "DefaultParagraphFont"
$headerChar.Priority = 99

$styleAnchorPara.Range.Style = "Footer"
$styleAnchorPara.Range.Style = $styleAnchorOriginal

$footerChar = $d.Styles.Add("FooterChar", 2)
$footerChar.NameLocal = "Footer Char"
$footerChar.BaseStyle = "DefaultParagraphFont"
$footerChar.Priority = 99

$headerStyle = $d.Styles("Header")
$headerStyle.NameLocal = "header"
$headerStyle.LinkStyle = "HeaderChar"
$headerStyle.UnhideWhenUsed = $true
$headerStyle.ParagraphFormat.SpaceAfter = 0
$headerStyle.ParagraphFormat.LineSpacingRule = 0

$footerStyle = $d.Styles("Footer")
$footerStyle.NameLocal = "footer"
$footerStyle.LinkStyle = "FooterChar"
$footerStyle.UnhideWhenUsed = $true
$footerStyle.ParagraphFormat.SpaceAfter = 0
$footerStyle.ParagraphFormat.LineSpacingRule = 0

$headerChar.LinkStyle = "Header"
$footerChar.LinkStyle = "Footer"

Write-Output "Removed completed item and refreshed document boilerplate."
